$wb = $excel.ActiveWorkbook

# ALC row 17
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 663195.75
$ws.Range("I17").Value = 400
$ws.Range("J17").Value = 695266.5
$ws.Range("K17").Value = 1200
$ws.Range("L17").Value = 2085799.5
$ws.Range("M17").Value = -1032
$ws.Range("N17").Value = -2086135.5

# ALC row 32
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3221.5
$ws.Range("J32").Value = 3221.5
$ws.Range("L32").Value = 3221.5
$ws.Range("N32").Value = -3873.5

# ALC row 109
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H109").Value = 17777.777
$ws.Range("J109").Value = 17777.777
$ws.Range("L109").Value = 17777.777
$ws.Range("N109").Value = -20551.777

# ALC row 137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 30304050
$ws.Range("I137").Value = 35714996
$ws.Range("J137").Value = 2763
$ws.Range("K137").Value = 107144988
$ws.Range("L137").Value = 8289
$ws.Range("M137").Value = -107142438
$ws.Range("N137").Value = -13389

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3947926.5
$ws.Range("I138").Value = 1264876.9
$ws.Range("J138").Value = 5557756.5
$ws.Range("K138").Value = 3794630.7
$ws.Range("L138").Value = 16673269.5
$ws.Range("M138").Value = -3789490.7
$ws.Range("N138").Value = -16683549.5

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1394.12
$ws.Range("I32").Value = 1348.2748
$ws.Range("J32").Value = 1857.6666
$ws.Range("K32").Value = 1348.2748
$ws.Range("L32").Value = 1857.6666
$ws.Range("M32").Value = -1061.2748
$ws.Range("N32").Value = -2431.6666

# ARM row 132
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2772.5454
$ws.Range("I132").Value = 2326.054
$ws.Range("J132").Value = 5132.5713
$ws.Range("K132").Value = 6978.162
$ws.Range("L132").Value = 15397.7139
$ws.Range("M132").Value = -4448.162
$ws.Range("N132").Value = -20457.7139

# ARM row 135
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H135").Value = 38177.4
$ws.Range("J135").Value = 38177.4
$ws.Range("L135").Value = 38177.4
$ws.Range("N135").Value = -48317.4

# ARM row 139
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H139").Value = 45380.445
$ws.Range("J139").Value = 45380.445
$ws.Range("L139").Value = 45380.445
$ws.Range("N139").Value = -55660.445

# BSM row 20
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2600.9
$ws.Range("I20").Value = 2333.3333
$ws.Range("J20").Value = 3002.25
$ws.Range("K20").Value = 2333.3333
$ws.Range("L20").Value = 3002.25
$ws.Range("M20").Value = -2086.3333
$ws.Range("N20").Value = -3496.25

# BSM row 135
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 38456.715
$ws.Range("J135").Value = 38456.715
$ws.Range("L135").Value = 38456.715
$ws.Range("N135").Value = -48596.715

# CRP row 121
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H121").Value = 40000
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 40000
$ws.Range("K121").Value = 0
$ws.Range("L121").ClearContents()
$ws.Range("N121").Value = -42620

# CRP row 132
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 3079.1516
$ws.Range("I132").Value = 2594.32
$ws.Range("K132").Value = 7782.960000000001
$ws.Range("M132").Value = -5252.960000000001

# CUL row 5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1302.4595
$ws.Range("I5").Value = 1024.84
$ws.Range("J5").Value = 1880.8334
$ws.Range("K5").Value = 3074.52
$ws.Range("L5").Value = 5642.5002
$ws.Range("M5").Value = -2962.52
$ws.Range("N5").Value = -5866.5002

# CUL row 112
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 1000000000
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 1000000000
$ws.Range("K112").Value = 0
$ws.Range("L112").ClearContents()
$ws.Range("N112").Value = -3000002216

# CUL row 122
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 928.7273
$ws.Range("I122").Value = 351.33334
$ws.Range("J122").Value = 1145.25
$ws.Range("K122").Value = 3162.00006
$ws.Range("L122").Value = 10307.25
$ws.Range("M122").Value = -712.0000600000003
$ws.Range("N122").Value = -15207.25

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 5954027.5
$ws.Range("J131").Value = 6946286
$ws.Range("L131").Value = 20838858
$ws.Range("N131").Value = -20848938

# CUL row 135
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 1302.4595
$ws.Range("I135").Value = 1024.84
$ws.Range("J135").Value = 1880.8334
$ws.Range("K135").Value = 9223.559999999999
$ws.Range("L135").Value = 16927.5006
$ws.Range("M135").Value = -6688.559999999999
$ws.Range("N135").Value = -21997.5006

# GSM row 70
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5915.121
$ws.Range("I70").Value = 6050.4346
$ws.Range("J70").Value = 5603.9
$ws.Range("K70").Value = 6050.4346
$ws.Range("L70").Value = 5603.9
$ws.Range("M70").Value = -5780.4346
$ws.Range("N70").Value = -6143.9

# GSM row 73
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 5915.121
$ws.Range("I73").Value = 6050.4346
$ws.Range("J73").Value = 5603.9
$ws.Range("K73").Value = 6050.4346
$ws.Range("L73").Value = 5603.9
$ws.Range("M73").Value = -5114.4346
$ws.Range("N73").Value = -7475.9

# GSM row 121
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H121").Value = 30000
$ws.Range("J121").Value = 30000
$ws.Range("L121").Value = 30000
$ws.Range("N121").Value = -33494

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1113179.9
$ws.Range("I122").Value = 1853473.9
$ws.Range("K122").Value = 5560421.699999999
$ws.Range("M122").Value = -5557971.699999999

# GSM row 126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1934.079
$ws.Range("I126").Value = 1470.2941
$ws.Range("J126").Value = 2309.524
$ws.Range("K126").Value = 4410.8823
$ws.Range("L126").Value = 6928.572
$ws.Range("M126").Value = -1940.8823
$ws.Range("N126").Value = -11868.572

# LTW row 2
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 1000
$ws.Range("J2").Value = 1000
$ws.Range("L2").Value = 1000
$ws.Range("N2").Value = -1224

# LTW row 55
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 542.4706
$ws.Range("I55").Value = 464.14285
$ws.Range("J55").Value = 597.3
$ws.Range("K55").Value = 464.14285
$ws.Range("L55").Value = 597.3
$ws.Range("M55").Value = -291.14285
$ws.Range("N55").Value = -943.3

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 3125
$ws.Range("I100").Value = 2800
$ws.Range("J100").Value = 3206.25
$ws.Range("K100").Value = 2800
$ws.Range("L100").Value = 3206.25
$ws.Range("M100").Value = -2259
$ws.Range("N100").Value = -4288.25

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 3227.0286
$ws.Range("I132").Value = 1786.4348
$ws.Range("J132").Value = 5988.1665
$ws.Range("K132").Value = 5359.3044
$ws.Range("L132").Value = 17964.4995
$ws.Range("M132").Value = -2829.3044
$ws.Range("N132").Value = -23024.4995

# WVR row 116
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H116").Value = 72680
$ws.Range("J116").Value = 72680
$ws.Range("L116").Value = 72680
$ws.Range("N116").Value = -81858

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2842.6826
$ws.Range("I132").Value = 3004.2766
$ws.Range("J132").Value = 2368
$ws.Range("K132").Value = 9012.8298
$ws.Range("L132").Value = 7104
$ws.Range("M132").Value = -6482.8298
$ws.Range("N132").Value = -12164
